# Apply updated symbol list values (price/volume/hour refresh) to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '244.32' }
    @{ Cell = "E2"; Value = '-1.10%' }
    @{ Cell = "G2"; Value = '15' }
    @{ Cell = "D3"; Value = '27.21' }
    @{ Cell = "E3"; Value = '3.55%' }
    @{ Cell = "G3"; Value = '15' }
    @{ Cell = "D4"; Value = '5.043' }
    @{ Cell = "E4"; Value = '-0.69%' }
    @{ Cell = "G4"; Value = '15' }
    @{ Cell = "D5"; Value = '0.05672' }
    @{ Cell = "E5"; Value = '0.99%' }
    @{ Cell = "G5"; Value = '15' }
    @{ Cell = "D6"; Value = '6.475' }
    @{ Cell = "E6"; Value = '-0.60%' }
    @{ Cell = "G6"; Value = '15' }
    @{ Cell = "D7"; Value = '0.8214' }
    @{ Cell = "E7"; Value = '1.01%' }
    @{ Cell = "G7"; Value = '15' }
    @{ Cell = "D8"; Value = '0.8442' }
    @{ Cell = "E8"; Value = '-0.62%' }
    @{ Cell = "G8"; Value = '15' }
    @{ Cell = "E9"; Value = '-1.02%' }
    @{ Cell = "G9"; Value = '15' }
    @{ Cell = "D10"; Value = '0.06923' }
    @{ Cell = "E10"; Value = '-0.60%' }
    @{ Cell = "G10"; Value = '15' }
    @{ Cell = "B11"; Value = 'BitrueCoin' }
    @{ Cell = "C11"; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = "D11"; Value = '0.02863' }
    @{ Cell = "E11"; Value = '1.36%' }
    @{ Cell = "G11"; Value = '15' }
    @{ Cell = "B12"; Value = 'BitMartToken' }
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = "D12"; Value = '0.09379' }
    @{ Cell = "E12"; Value = '-0.24%' }
    @{ Cell = "G12"; Value = '15' }
    @{ Cell = "B13"; Value = 'BitForexToken' }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = "D13"; Value = '0.001526' }
    @{ Cell = "E13"; Value = '0.88%' }
    @{ Cell = "G13"; Value = '15' }
    @{ Cell = "B14"; Value = 'CoinExToken' }
    @{ Cell = "C14"; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' }
    @{ Cell = "D14"; Value = '0.04116' }
    @{ Cell = "E14"; Value = '-12.28%' }
    @{ Cell = "G14"; Value = '15' }
    @{ Cell = "B15"; Value = 'One' }
    @{ Cell = "C15"; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one' }
    @{ Cell = "D15"; Value = '0.0006010' }
    @{ Cell = "E15"; Value = '0.85%' }
    @{ Cell = "G15"; Value = '15' }
    @{ Cell = "B16"; Value = 'TigerCash' }
    @{ Cell = "C16"; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = "D16"; Value = '0.006182' }
    @{ Cell = "E16"; Value = '0.83%' }
    @{ Cell = "G16"; Value = '15' }
    @{ Cell = "B17"; Value = 'LEO' }
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = "D17"; Value = '3.510' }
    @{ Cell = "E17"; Value = '-1.89%' }
    @{ Cell = "G17"; Value = '15' }
    @{ Cell = "B18"; Value = 'GateToken' }
    @{ Cell = "C18"; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = "D18"; Value = '3.000' }
    @{ Cell = "E18"; Value = '-1.89%' }
    @{ Cell = "G18"; Value = '15' }
    @{ Cell = "B19"; Value = 'BTSEToken' }
    @{ Cell = "C19"; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' }
    @{ Cell = "D19"; Value = '2.310' }
    @{ Cell = "E19"; Value = '9.04%' }
    @{ Cell = "G19"; Value = '15' }
    @{ Cell = "B20"; Value = 'BitpandaEcosystemToken' }
    @{ Cell = "C20"; Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best' }
    @{ Cell = "D20"; Value = '0.3113' }
    @{ Cell = "E20"; Value = '-2.14%' }
    @{ Cell = "G20"; Value = '15' }
    @{ Cell = "B21"; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = "C21"; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = "D21"; Value = '0.03168' }
    @{ Cell = "E21"; Value = '0.41%' }
    @{ Cell = "G21"; Value = '15' }
    @{ Cell = "D22"; Value = '0.1254' }
    @{ Cell = "E22"; Value = '-4.99%' }
    @{ Cell = "G22"; Value = '15' }
    @{ Cell = "D23"; Value = '3.554' }
    @{ Cell = "E23"; Value = '-5.06%' }
    @{ Cell = "G23"; Value = '15' }
    @{ Cell = "E24"; Value = '-0.12%' }
    @{ Cell = "G24"; Value = '15' }
    @{ Cell = "D25"; Value = '0.001221' }
    @{ Cell = "E25"; Value = '-2.25%' }
    @{ Cell = "G25"; Value = '15' }
    @{ Cell = "D26"; Value = '0.003870' }
    @{ Cell = "E26"; Value = '-16.19%' }
    @{ Cell = "G26"; Value = '15' }
    @{ Cell = "E27"; Value = '2.09%' }
    @{ Cell = "G27"; Value = '15' }
    @{ Cell = "E28"; Value = '-25.76%' }
    @{ Cell = "G28"; Value = '15' }
    @{ Cell = "G29"; Value = '15' }
    @{ Cell = "G30"; Value = '15' }
    @{ Cell = "G31"; Value = '15' }
    @{ Cell = "G32"; Value = '15' }
    @{ Cell = "G33"; Value = '15' }
    @{ Cell = "G34"; Value = '15' }
    @{ Cell = "G35"; Value = '15' }
    @{ Cell = "G36"; Value = '15' }
    @{ Cell = "G37"; Value = '15' }
    @{ Cell = "G38"; Value = '15' }
    @{ Cell = "G39"; Value = '15' }
    @{ Cell = "D40"; Value = '0.03665' }
    @{ Cell = "E40"; Value = '-0.16%' }
    @{ Cell = "G40"; Value = '15' }
    @{ Cell = "B41"; Value = 'KickToken' }
    @{ Cell = "C41"; Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick' }
    @{ Cell = "D41"; Value = '0.006050' }
    @{ Cell = "E41"; Value = '75.66%' }
    @{ Cell = "G41"; Value = '15' }
    @{ Cell = "B42"; Value = 'BKEXToken' }
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk' }
    @{ Cell = "D42"; Value = '0.1052' }
    @{ Cell = "E42"; Value = '-23.09%' }
    @{ Cell = "G42"; Value = '15' }
    @{ Cell = "D43"; Value = '0.002293' }
    @{ Cell = "E43"; Value = '-13.80%' }
    @{ Cell = "G43"; Value = '15' }
    @{ Cell = "D44"; Value = '0.009537' }
    @{ Cell = "E44"; Value = '11.05%' }
    @{ Cell = "G44"; Value = '15' }
    @{ Cell = "D45"; Value = '0.00005313' }
    @{ Cell = "E45"; Value = '0.38%' }
    @{ Cell = "G45"; Value = '15' }
    @{ Cell = "E46"; Value = '0.01%' }
    @{ Cell = "G46"; Value = '15' }
    @{ Cell = "E47"; Value = '-15.41%' }
    @{ Cell = "G47"; Value = '15' }
    @{ Cell = "D48"; Value = '0.002568' }
    @{ Cell = "E48"; Value = '24.09%' }
    @{ Cell = "G48"; Value = '15' }
    @{ Cell = "E49"; Value = '0.01%' }
    @{ Cell = "G49"; Value = '15' }
    @{ Cell = "E50"; Value = '0.01%' }
    @{ Cell = "G50"; Value = '15' }
    @{ Cell = "G51"; Value = '15' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (prices, percentages, hours)
    # are not silently converted to numbers by Excel.
    $range.NumberFormat = '@'
    $range.Value = $u.Value
}

